$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CIList")

# Append two more rows duplicating existing CI values (SE_AAA_AAA04SE, SE_AAA_AAA98SE)
$ws.Range("A10").Value = "SE_AAA_AAA04SE"
$ws.Range("A11").Value = "SE_AAA_AAA98SE"

# Update selection to match the new active cell after data entry
$ws.Range("A11").Select()
